# Update cryptocurrency price/volume data to the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.008.22"
$ws.Range("E2").Value = "  +3.04%  "
$ws.Range("D3").Value = "3.403.39"
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'584.28"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "'181.11"
$ws.Range("E6").Value = "  +1.92%  "
$ws.Range("D7").Value = "'0.600"
$ws.Range("E7").Value = "  +1.28%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +9.81%  "
$ws.Range("D10").Value = "'0.596"
$ws.Range("E10").Value = "  +2.17%  "
$ws.Range("D11").Value = "'48.81"
$ws.Range("E12").Value = "  +5.13%  "
$ws.Range("D13").Value = "'689.81"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("D14").Value = "'8.71"
$ws.Range("E14").Value = "  +3.02%  "
$ws.Range("D15").Value = "3.955.29"
$ws.Range("E15").Value = "  +2.06%  "
$ws.Range("D16").Value = "69.928.29"
$ws.Range("E16").Value = "  +2.85%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.413.42"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.121"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("D19").Value = "'17.76"
$ws.Range("E20").Value = "  +2.04%  "
$ws.Range("D21").Value = "'0.921"
$ws.Range("E21").Value = "  +2.87%  "
$ws.Range("D22").Value = "'17.32"
$ws.Range("E22").Value = "  +2.22%  "
$ws.Range("D23").Value = "'5.38"
$ws.Range("E23").Value = "  -0.81%  "
$ws.Range("D24").Value = "'102.77"
$ws.Range("E24").Value = "  +2.43%  "
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("E26").Value = "  +1.14%  "
$ws.Range("D27").Value = "'9.66"
$ws.Range("E27").Value = "  +1.88%  "
$ws.Range("D28").Value = "'33.85"
$ws.Range("E28").Value = "  +2.34%  "
$ws.Range("D29").Value = "'8.84"
$ws.Range("E29").Value = "  +3.15%  "
$ws.Range("D30").Value = "'6.99"
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("D31").Value = "'3.72"
$ws.Range("E31").Value = "  +10.65%  "
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("D33").Value = "'556.74"
$ws.Range("E33").Value = "  -3.71%  "
$ws.Range("D35").Value = "'58.67"
$ws.Range("E35").Value = "  +2.26%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "3.668.37"
$ws.Range("E37").Value = "  -1.99%  "
$ws.Range("D38").Value = "'0.140"
$ws.Range("E38").Value = "  +4.32%  "
$ws.Range("D39").Value = "'35.65"
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("E40").Value = "  +8.96%  "
$ws.Range("D41").Value = "'3.34"
$ws.Range("E41").Value = "  +5.22%  "
$ws.Range("E42").Value = "  +3.26%  "
$ws.Range("E43").Value = "  +4.56%  "
$ws.Range("E44").Value = "  +1.43%  "
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "'2.68"
$ws.Range("E45").Value = "  +2.06%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.130"
$ws.Range("E46").Value = "  +0.87%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'1.39"
$ws.Range("E47").Value = "  +4.41%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'130.39"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("B50").Value = "CoreDAO"
$ws.Range("C50").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D50").Value = "'2.62"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "'7.50"
$ws.Range("E51").Value = "  +2.06%  "
